$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") for rows 2-9 is bumped by one day (2023-10-08 -> 2023-10-09),
# i.e. serial date value 45207 -> 45208.
$ws.Range("C2:C9").Value = 45208
